# Edit values of capacitors near tone generator + updated BOM
# Adds 5 new capacitor BOM rows (1800pF..3900pF, Farnell) plus a subtotal.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New BOM rows 40-44: additional tone-generator capacitors ---

# Row 40 - 1800pF
$ws.Range("A40").Value = "1800pF"
$ws.Range("B40").Value = "Farnell"
$ws.Range("C40").Value = 2524826
$ws.Range("D40").Value = 0.248
$ws.Range("E40").Value = 1

# Row 41 - 2200pF
$ws.Range("A41").Value = "2200pF"
$ws.Range("B41").Value = "Farnell"
$ws.Range("C41").Value = 2496845
$ws.Range("D41").Value = 0.185
$ws.Range("E41").Value = 1

# Row 42 - 2700pF
$ws.Range("A42").Value = "2700pF"
$ws.Range("B42").Value = "Farnell"
$ws.Range("C42").Value = 2524833
$ws.Range("D42").Value = 0.276
$ws.Range("E42").Value = 1

# Row 43 - 3300pF
$ws.Range("A43").Value = "3300pF"
$ws.Range("B43").Value = "Farnell"
$ws.Range("C43").Value = 2496853
$ws.Range("D43").Value = 0.182
$ws.Range("E43").Value = 1

# Row 44 - 3900pF
$ws.Range("A44").Value = "3900pF"
$ws.Range("B44").Value = "Farnell"
$ws.Range("C44").Value = 2524841
$ws.Range("D44").Value = 0.248
$ws.Range("E44").Value = 1

# Row 47 - subtotal of the new capacitors
$ws.Range("D47").Formula = "=SUM(D40:D44)"

# Update the visible selection to match the new bottom of the sheet
$ws.Range("D48").Select()
